$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.358.49'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '1.791.40'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.556'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.87'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.79%  '
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0690'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0946'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '2.053.17'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('D14').Value = '1.795.22'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '34.384.65'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '245.72'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.71%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '169.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('E26').Value = '  +3.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.57'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0527'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.24%  '
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('D35').Value = '1.421.27'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.688'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.57'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.46%  '
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0191'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '84.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.947'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.69%  '
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('E43').Value = '  +1.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0525'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('D48').Value = '1.955.22'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.45'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.61%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -1.18%  '
